# README.xlsx: correct folder names on the README tab.
#
# The README sheet has two summary lines describing the other two tabs.
# Previously they were mislabeled/swapped ("Input_files" described the
# ModelInput tab but used the word "directory", etc). Fix them so:
#   A13 -> bold "ModelOutput" lead-in, describing the 'ModelOutput' folder
#   A14 -> bold "ModelInput" lead-in, describing the 'ModelInput' folder
# and both now say "folder" instead of "directory".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("README")

# --- A13: was "Input_files: lists all files contained in the directory
#            'ModelInput'. "
#     now  "ModelOutput: lists model output files contained in the folder
#            'ModelOutput'. "
$a13 = $ws.Range("A13")
$a13.Value = "ModelOutput: lists model output files contained in the folder 'ModelOutput'. "
$a13Bold = $a13.Characters(1, 11)
$a13Bold.Font.Bold = $true
$a13Bold.Font.Size = 12
$a13Rest = $a13.Characters(12, 67)
$a13Rest.Font.Bold = $false
$a13Rest.Font.Size = 12

# --- A14: was "Output_files: lists model output files contained in the
#            directory 'ModelOutput'. "
#     now  "ModelInput: lists  model input files contained in the folder
#            'ModelInput'. "
$a14 = $ws.Range("A14")
$a14.Value = "ModelInput: lists  model input files contained in the folder 'ModelInput'. "
$a14Bold = $a14.Characters(1, 10)
$a14Bold.Font.Bold = $true
$a14Bold.Font.Size = 12
$a14Rest = $a14.Characters(11, 65)
$a14Rest.Font.Bold = $false
$a14Rest.Font.Size = 12

# The saved view previously had A13 highlighted (left over from editing);
# reset the active selection back to the top of the sheet.
[void]$ws.Range("A1").Select()
